# Auto-generated update script for resum_diari_meteocat.xlsx
# Commit: Update automàtic: dades i banners [2026-02-12 23:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-12 23:48:50"
$ws.Range("K2").Value = "7.0 MJ/m2"
$ws.Range("N2").Value = "-0.1 °C 23:19 TU"
$ws.Range("E3").Value = "2026-02-12 23:48:53"
$ws.Range("E4").Value = "2026-02-12 23:48:55"
$ws.Range("J4").Value = "1000.1 hPa"
$ws.Range("E5").Value = "2026-02-12 23:48:58"
$ws.Range("E6").Value = "2026-02-12 23:49:01"
$ws.Range("J6").Value = "1000.0 hPa"
$ws.Range("N6").Value = "10.8 °C 23:29 TU"
$ws.Range("O6").Value = "15.5 °C"
$ws.Range("E7").Value = "2026-02-12 23:49:04"
$ws.Range("J7").Value = "1002.5 hPa"
$ws.Range("O7").Value = "16.9 °C"
$ws.Range("E8").Value = "2026-02-12 23:49:06"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "46%"
$ws.Range("J8").Value = "1001.9 hPa"
$ws.Range("O8").Value = "12.9 °C"
$ws.Range("E9").Value = "2026-02-12 23:49:09"
$ws.Range("O9").Value = "12.9 °C"
$ws.Range("E10").Value = "2026-02-12 23:49:12"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "49%"
$ws.Range("N10").Value = "6.7 °C 23:29 TU"
$ws.Range("O10").Value = "14.3 °C"
$ws.Range("E11").Value = "2026-02-12 23:49:15"
$ws.Range("N11").Value = "2.7 °C 23:29 TU"
$ws.Range("O11").Value = "8.8 °C"
$ws.Range("E12").Value = "2026-02-12 23:49:17"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "69%"
$ws.Range("E13").Value = "2026-02-12 23:49:20"
$ws.Range("J13").Value = "1002.7 hPa"
$ws.Range("K13").Value = "13.3 MJ/m2"
$ws.Range("N13").Value = "2.9 °C 23:25 TU"
$ws.Range("O13").Value = "7.3 °C"
$ws.Range("E14").Value = "2026-02-12 23:49:23"
$ws.Range("E15").Value = "2026-02-12 23:49:25"
$ws.Range("E16").Value = "2026-02-12 23:49:28"
$ws.Range("E17").Value = "2026-02-12 23:49:31"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "73%"
$ws.Range("E18").Value = "2026-02-12 23:49:33"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "39%"
$ws.Range("J18").Value = "1000.4 hPa"
$ws.Range("O18").Value = "16.1 °C"
$ws.Range("E19").Value = "2026-02-12 23:49:36"
$ws.Range("N19").Value = "4.7 °C 23:04 TU"
$ws.Range("E20").Value = "2026-02-12 23:49:39"
$ws.Range("E21").Value = "2026-02-12 23:49:42"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "50%"
$ws.Range("J21").Value = "1003.2 hPa"
$ws.Range("N21").Value = "4.4 °C 23:28 TU"
$ws.Range("O21").Value = "8.8 °C"
$ws.Range("E22").Value = "2026-02-12 23:49:45"
$ws.Range("E23").Value = "2026-02-12 23:49:47"
$ws.Range("K23").Value = "10.8 MJ/m2"
$ws.Range("E24").Value = "2026-02-12 23:49:50"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "61%"
$ws.Range("J24").Value = "1007.0 hPa"
$ws.Range("O24").Value = "11.4 °C"
$ws.Range("E25").Value = "2026-02-12 23:49:53"
$ws.Range("E26").Value = "2026-02-12 23:49:55"
$ws.Range("J26").Value = "999.9 hPa"
$ws.Range("O26").Value = "5.6 °C"
$ws.Range("E27").Value = "2026-02-12 23:49:58"
$ws.Range("K27").Value = "13.7 MJ/m2"
$ws.Range("E28").Value = "2026-02-12 23:50:01"
$ws.Range("J28").Value = "999.9 hPa"
$ws.Range("O28").Value = "13.4 °C"
$ws.Range("E29").Value = "2026-02-12 23:50:04"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "63%"
$ws.Range("O29").Value = "13.5 °C"
$ws.Range("E30").Value = "2026-02-12 23:50:06"
$ws.Range("J30").Value = "1000.2 hPa"
$ws.Range("O30").Value = "11.6 °C"
$ws.Range("E31").Value = "2026-02-12 23:50:09"
$ws.Range("J31").Value = "999.6 hPa"
$ws.Range("E32").Value = "2026-02-12 23:50:11"
$ws.Range("O32").Value = "7.9 °C"
$ws.Range("E33").Value = "2026-02-12 23:50:14"
$ws.Range("J33").Value = "1002.4 hPa"
$ws.Range("N33").Value = "1.9 °C 23:29 TU"
$ws.Range("O33").Value = "6.3 °C"
$ws.Range("E34").Value = "2026-02-12 23:50:17"
$ws.Range("E35").Value = "2026-02-12 23:50:20"
$ws.Range("E36").Value = "2026-02-12 23:50:22"
$ws.Range("J36").Value = "1000.5 hPa"
$ws.Range("O36").Value = "14.2 °C"
$ws.Range("E37").Value = "2026-02-12 23:50:25"
$ws.Range("E38").Value = "2026-02-12 23:50:28"
$ws.Range("O38").Value = "15.6 °C"
$ws.Range("E39").Value = "2026-02-12 23:50:30"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "63%"
$ws.Range("E40").Value = "2026-02-12 23:50:33"
$ws.Range("J40").Value = "1004.0 hPa"
$ws.Range("N40").Value = "2.5 °C 23:25 TU"
$ws.Range("O40").Value = "8.9 °C"
$ws.Range("E41").Value = "2026-02-12 23:50:36"
$ws.Range("O41").Value = "17.0 °C"
$ws.Range("E42").Value = "2026-02-12 23:50:38"
$ws.Range("O42").Value = "13.4 °C"
$ws.Range("E43").Value = "2026-02-12 23:50:41"
$ws.Range("E44").Value = "2026-02-12 23:50:44"
$ws.Range("N44").Value = "-6.2 °C 23:14 TU"
$ws.Range("O44").Value = "-3.2 °C"
$ws.Range("E45").Value = "2026-02-12 23:50:46"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "56%"
$ws.Range("J45").Value = "1005.7 hPa"
$ws.Range("N45").Value = "0.8 °C 23:29 TU"
$ws.Range("O45").Value = "6.5 °C"
$ws.Range("E46").Value = "2026-02-12 23:50:49"
$ws.Range("N46").Value = "10.8 °C 23:29 TU"
$ws.Range("O46").Value = "15.5 °C"
